# 自动更新Excel文件 - 2025-10-09 23:11:46
# 每日巡检：剩余天数递减，若已到期（剩余=1）则按总天数续约并刷新开始时间

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$renewDate = 20251010

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $startDate = $ws.Cells.Item($r, 6).Value2

    if ($remaining -eq $null -or $total -eq $null) {
        continue
    }

    # 开始时间必须是合法的 8 位 yyyymmdd 日期，否则跳过该行（数据异常）
    $dateText = [string]$startDate
    if ($dateText.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        # 到期续约：剩余天数重置为总天数，开始时间更新为今天
        $ws.Cells.Item($r, 5).Value = $total
        $ws.Cells.Item($r, 6).Value = $renewDate
    } else {
        # 正常巡检：剩余天数减一
        $ws.Cells.Item($r, 5).Value = $remaining - 1
    }
}
